$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3334.1667
$ws.Range("I62").Value = 3001.25
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 3001.25
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -2377.25
$ws.Range("N62").Value = -5248

$ws.Range("H64").Value = 3983.3333
$ws.Range("I64").Value = 3487.2856
$ws.Range("J64").Value = 4156.95
$ws.Range("K64").Value = 3487.2856
$ws.Range("L64").Value = 4156.95
$ws.Range("M64").Value = -3239.2856
$ws.Range("N64").Value = -4652.95

$ws.Range("H65").Value = 3334.1667
$ws.Range("I65").Value = 3001.25
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 15006.25
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -11886.25
$ws.Range("N65").Value = -26240

$ws.Range("H67").Value = 3983.3333
$ws.Range("I67").Value = 3487.2856
$ws.Range("J67").Value = 4156.95
$ws.Range("K67").Value = 3487.2856
$ws.Range("L67").Value = 4156.95
$ws.Range("M67").Value = -2629.2856
$ws.Range("N67").Value = -5872.95

$ws.Range("H132").Value = 8175961
$ws.Range("I132").Value = 11153.4
$ws.Range("K132").Value = 33460.2
$ws.Range("M132").Value = -30930.2

$ws.Range("H137").Value = 3849569.2
$ws.Range("I137").Value = 7145457
$ws.Range("J137").Value = 4366.4165
$ws.Range("K137").Value = 21436371
$ws.Range("L137").Value = 13099.2495
$ws.Range("M137").Value = -21433821
$ws.Range("N137").Value = -18199.2495

$ws.Range("H138").Value = 1874763.1
$ws.Range("I138").Value = 1320.5555
$ws.Range("J138").Value = 2349720.5
$ws.Range("K138").Value = 3961.6665
$ws.Range("L138").Value = 7049161.5
$ws.Range("M138").Value = 1178.3335
$ws.Range("N138").Value = -7059441.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 1004
$ws.Range("I10").Value = 1004
$ws.Range("K10").Value = 1004
$ws.Range("M10").Value = -834

$ws.Range("H122").Value = 7409822
$ws.Range("I122").Value = 2927.3
$ws.Range("J122").Value = 22223610
$ws.Range("K122").Value = 8781.900000000001
$ws.Range("L122").Value = 66670830
$ws.Range("M122").Value = -6331.900000000001
$ws.Range("N122").Value = -66675730

$ws.Range("H123").Value = 58333.332
$ws.Range("J123").Value = 58333.332
$ws.Range("L123").Value = 58333.332
$ws.Range("N123").Value = -68133.33199999999

$ws.Range("H124").Value = 22554.834
$ws.Range("J124").Value = 22554.834
$ws.Range("L124").Value = 22554.834
$ws.Range("N124").Value = -32374.834

$ws.Range("H131").Value = 55789.473
$ws.Range("J131").Value = 55789.473
$ws.Range("L131").Value = 55789.473
$ws.Range("N131").Value = -65869.473

$ws.Range("H132").Value = 252481.25
$ws.Range("I132").Value = 252359
$ws.Range("J132").Value = 252603.5
$ws.Range("K132").Value = 757077
$ws.Range("L132").Value = 757810.5
$ws.Range("M132").Value = -754547
$ws.Range("N132").Value = -762870.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1575
$ws.Range("I99").Value = 1572
$ws.Range("J99").Value = 1590
$ws.Range("K99").Value = 1572
$ws.Range("L99").Value = 1590
$ws.Range("M99").Value = -74
$ws.Range("N99").Value = -4586

$ws.Range("H137").Value = 47830
$ws.Range("J137").Value = 47830
$ws.Range("L137").Value = 47830
$ws.Range("N137").Value = -58030

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2768.0303
$ws.Range("I31").Value = 1340.6957
$ws.Range("J31").Value = 3199.9868
$ws.Range("K31").Value = 1340.6957
$ws.Range("L31").Value = 3199.9868
$ws.Range("M31").Value = -1045.6957
$ws.Range("N31").Value = -3789.9868

$ws.Range("H34").Value = 2768.0303
$ws.Range("I34").Value = 1340.6957
$ws.Range("J34").Value = 3199.9868
$ws.Range("K34").Value = 1340.6957
$ws.Range("L34").Value = 3199.9868
$ws.Range("M34").Value = -1138.6957
$ws.Range("N34").Value = -3603.9868

$ws.Range("H58").Value = 25644030
$ws.Range("I58").Value = 32261396
$ws.Range("J58").Value = 1737
$ws.Range("K58").Value = 32261396
$ws.Range("L58").Value = 1737
$ws.Range("M58").Value = -32261193
$ws.Range("N58").Value = -2143

$ws.Range("H136").Value = 25644030
$ws.Range("I136").Value = 32261396
$ws.Range("J136").Value = 1737
$ws.Range("K136").Value = 96784188
$ws.Range("L136").Value = 5211
$ws.Range("M136").Value = -96781638
$ws.Range("N136").Value = -10311

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 3750
$ws.Range("J49").Value = 3750
$ws.Range("L49").Value = 11250
$ws.Range("N49").Value = -11562

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1366910.2
$ws.Range("I11").Value = 2000130.2
$ws.Range("J11").Value = 100470
$ws.Range("K11").Value = 2000130.2
$ws.Range("L11").Value = 100470
$ws.Range("M11").Value = -1999991.2
$ws.Range("N11").Value = -100748

$ws.Range("H99").Value = 6265
$ws.Range("I99").Value = 6265
$ws.Range("K99").Value = 6265
$ws.Range("M99").Value = -4019

$ws.Range("H102").Value = 1071.1428
$ws.Range("I102").Value = 1059.6
$ws.Range("J102").Value = 1100
$ws.Range("K102").Value = 1059.6
$ws.Range("L102").Value = 1100
$ws.Range("M102").Value = 562.4000000000001
$ws.Range("N102").Value = -4344

$ws.Range("H122").Value = 4835.1816
$ws.Range("I122").Value = 4748
$ws.Range("J122").Value = 4939.8
$ws.Range("K122").Value = 14244
$ws.Range("L122").Value = 14819.4
$ws.Range("M122").Value = -11794
$ws.Range("N122").Value = -19719.4

$ws.Range("H132").Value = 104086.4
$ws.Range("I132").Value = 104359.4
$ws.Range("K132").Value = 313078.2
$ws.Range("M132").Value = -310548.2

$ws.Range("H136").Value = 26800
$ws.Range("J136").Value = 26800
$ws.Range("L136").Value = 80400
$ws.Range("N136").Value = -85500

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1485.2273
$ws.Range("I16").Value = 621.94116
$ws.Range("J16").Value = 4420.4
$ws.Range("K16").Value = 621.94116
$ws.Range("L16").Value = 4420.4
$ws.Range("M16").Value = -451.94116
$ws.Range("N16").Value = -4760.4

$ws.Range("H22").Value = 774.8261
$ws.Range("I22").Value = 408.27274
$ws.Range("J22").Value = 1110.8334
$ws.Range("K22").Value = 408.27274
$ws.Range("L22").Value = 1110.8334
$ws.Range("M22").Value = -113.27274
$ws.Range("N22").Value = -1700.8334

$ws.Range("H27").Value = 774.8261
$ws.Range("I27").Value = 408.27274
$ws.Range("J27").Value = 1110.8334
$ws.Range("K27").Value = 408.27274
$ws.Range("L27").Value = 1110.8334
$ws.Range("M27").Value = -301.27274
$ws.Range("N27").Value = -1324.8334

$ws.Range("H99").Value = 11629.5
$ws.Range("I99").Value = 11629.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 11629.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -8634.5
$ws.Range("N99").ClearContents()

$ws.Range("H100").Value = 35952.863
$ws.Range("I100").Value = 92136.37
$ws.Range("J100").Value = 1618.5
$ws.Range("K100").Value = 92136.37
$ws.Range("L100").Value = 1618.5
$ws.Range("M100").Value = -91595.37
$ws.Range("N100").Value = -2700.5

$ws.Range("H122").Value = 3596.5833
$ws.Range("I122").Value = 3219.875
$ws.Range("J122").Value = 4350
$ws.Range("K122").Value = 9659.625
$ws.Range("L122").Value = 13050
$ws.Range("M122").Value = -7209.625
$ws.Range("N122").Value = -17950

$ws.Range("H132").Value = 66518.75
$ws.Range("I132").Value = 3437.75
$ws.Range("J132").Value = 129599.75
$ws.Range("K132").Value = 10313.25
$ws.Range("L132").Value = 388799.25
$ws.Range("M132").Value = -7783.25
$ws.Range("N132").Value = -393859.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 56666.668
$ws.Range("J131").Value = 56666.668
$ws.Range("L131").Value = 56666.668
$ws.Range("N131").Value = -66746.66800000001

$ws.Range("H132").Value = 253709.62
$ws.Range("I132").Value = 253669.75
$ws.Range("K132").Value = 761009.25
$ws.Range("M132").Value = -758479.25

$ws.Range("H136").Value = 108364.734
$ws.Range("I136").Value = 87442.836
$ws.Range("K136").Value = 262328.508
$ws.Range("M136").Value = -259778.508
